$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") "318.41"
Set-TextValue $ws.Range("E2") "4.25%"
Set-TextValue $ws.Range("D3") "39.50"
Set-TextValue $ws.Range("E3") "3.34%"
Set-TextValue $ws.Range("D4") "5.117"
Set-TextValue $ws.Range("E4") "0.70%"
Set-TextValue $ws.Range("D5") "0.08206"
Set-TextValue $ws.Range("E5") "1.86%"
Set-TextValue $ws.Range("D6") "2.044"
Set-TextValue $ws.Range("E6") "5.52%"
Set-TextValue $ws.Range("D7") "8.271"
Set-TextValue $ws.Range("E7") "4.06%"
Set-TextValue $ws.Range("D8") "0.9333"
Set-TextValue $ws.Range("E8") "0.39%"
Set-TextValue $ws.Range("D9") "0.1412"
Set-TextValue $ws.Range("E9") "-3.67%"
Set-TextValue $ws.Range("D10") "0.1987"
Set-TextValue $ws.Range("E10") "3.44%"
Set-TextValue $ws.Range("D11") "0.09108"
Set-TextValue $ws.Range("E11") "1.33%"
Set-TextValue $ws.Range("D12") "0.03559"
Set-TextValue $ws.Range("E12") "1.20%"
Set-TextValue $ws.Range("D13") "0.09816"
Set-TextValue $ws.Range("E13") "0.28%"
Set-TextValue $ws.Range("D14") "0.001401"
Set-TextValue $ws.Range("E14") "0.32%"
Set-TextValue $ws.Range("D15") "0.006328"
Set-TextValue $ws.Range("E15") "4.32%"
Set-TextValue $ws.Range("D16") "3.658"
Set-TextValue $ws.Range("E16") "-1.69%"
Set-TextValue $ws.Range("D17") "4.278"
Set-TextValue $ws.Range("E17") "2.19%"
Set-TextValue $ws.Range("E18") "-3.38%"
Set-TextValue $ws.Range("D19") "0.3460"
Set-TextValue $ws.Range("E19") "-0.05%"
Set-TextValue $ws.Range("D20") "0.1303"
Set-TextValue $ws.Range("E20") "-0.69%"
Set-TextValue $ws.Range("D21") "4.904"
Set-TextValue $ws.Range("E21") "1.99%"
Set-TextValue $ws.Range("D22") "0.2450"
Set-TextValue $ws.Range("E22") "1.81%"
Set-TextValue $ws.Range("D23") "0.04330"
Set-TextValue $ws.Range("E23") "-1.25%"
Set-TextValue $ws.Range("D24") "0.001226"
Set-TextValue $ws.Range("E24") "-0.55%"
Set-TextValue $ws.Range("D25") "0.004783"
Set-TextValue $ws.Range("E25") "12.58%"
Set-TextValue $ws.Range("D26") "0.0001301"
Set-TextValue $ws.Range("E26") "-0.11%"
Set-TextValue $ws.Range("E27") "-10.02%"
Set-TextValue $ws.Range("D39") "0.02231"
Set-TextValue $ws.Range("E39") "8.19%"
Set-TextValue $ws.Range("D40") "0.05252"
Set-TextValue $ws.Range("E40") "4.44%"
Set-TextValue $ws.Range("D41") "0.007505"
Set-TextValue $ws.Range("E41") "0.39%"
Set-TextValue $ws.Range("D42") "0.009792"
Set-TextValue $ws.Range("D43") "0.1380"
Set-TextValue $ws.Range("E43") "2.50%"
Set-TextValue $ws.Range("D44") "0.002122"
Set-TextValue $ws.Range("E44") "-1.02%"
Set-TextValue $ws.Range("D45") "0.009517"
Set-TextValue $ws.Range("E45") "6.07%"
Set-TextValue $ws.Range("D46") "0.00006450"
Set-TextValue $ws.Range("E46") "4.19%"
Set-TextValue $ws.Range("E47") "-0.13%"
Set-TextValue $ws.Range("B48") "BOLO"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws.Range("D48") "0.002770"
Set-TextValue $ws.Range("E48") "-1.47%"
Set-TextValue $ws.Range("B49") "CoinbaseStockToken"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws.Range("D49") "0.001200"
Set-TextValue $ws.Range("E49") "-25.10%"
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "-0.13%"
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "-0.13%"
